# Auto-generated script applying scheduled market-price/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 4000
$ws.Range("I34").Value = 944.44446
$ws.Range("J34").Value = 13166.667
$ws.Range("K34").Value = 944.44446
$ws.Range("L34").Value = 13166.667
$ws.Range("M34").Value = -741.44446
$ws.Range("N34").Value = -13572.667
$ws.Range("H36").Value = 4000
$ws.Range("I36").Value = 944.44446
$ws.Range("J36").Value = 13166.667
$ws.Range("K36").Value = 944.44446
$ws.Range("L36").Value = 13166.667
$ws.Range("M36").Value = -229.44446
$ws.Range("N36").Value = -14596.667
$ws.Range("H48").Value = 3000
$ws.Range("J48").Value = 3000
$ws.Range("L48").Value = 9000
$ws.Range("N48").Value = -9584
$ws.Range("H56").Value = 3000
$ws.Range("J56").Value = 3000
$ws.Range("L56").Value = 9000
$ws.Range("N56").Value = -10068
$ws.Range("H99").Value = 2382.4
$ws.Range("I99").Value = 2778
$ws.Range("J99").Value = 800
$ws.Range("K99").Value = 8334
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -6836
$ws.Range("N99").Value = -5396
$ws.Range("H101").Value = 14769.857
$ws.Range("I101").Value = 14769.857
$ws.Range("K101").Value = 44309.571
$ws.Range("M101").Value = -42687.571
$ws.Range("H112").Value = 967.9666999999999
$ws.Range("I112").Value = 550
$ws.Range("J112").Value = 982.37933
$ws.Range("K112").Value = 1650
$ws.Range("L112").Value = 2947.13799
$ws.Range("M112").Value = -542
$ws.Range("N112").Value = -5163.13799
$ws.Range("H129").Value = 862.8261
$ws.Range("J129").Value = 902.75
$ws.Range("L129").Value = 2708.25
$ws.Range("N129").Value = -12708.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 2603
$ws.Range("I39").Value = 2603
$ws.Range("K39").Value = 2603
$ws.Range("M39").Value = -2083
$ws.Range("H76").Value = 32644
$ws.Range("J76").Value = 32644
$ws.Range("L76").Value = 32644
$ws.Range("N76").Value = -33320
$ws.Range("H79").Value = 32644
$ws.Range("J79").Value = 32644
$ws.Range("L79").Value = 32644
$ws.Range("N79").Value = -34984
$ws.Range("H122").Value = 2436.2222
$ws.Range("I122").Value = 2597.889
$ws.Range("J122").Value = 2274.5557
$ws.Range("K122").Value = 7793.667
$ws.Range("L122").Value = 6823.6671
$ws.Range("M122").Value = -5343.667
$ws.Range("N122").Value = -11723.6671
$ws.Range("H131").Value = 48708.75
$ws.Range("J131").Value = 48708.75
$ws.Range("L131").Value = 48708.75
$ws.Range("N131").Value = -58788.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 183828.9
$ws.Range("I105").Value = 202056
$ws.Range("K105").Value = 202056
$ws.Range("M105").Value = -200309
$ws.Range("H134").Value = 2377.742
$ws.Range("I134").Value = 2386.5
$ws.Range("K134").Value = 7159.5
$ws.Range("M134").Value = -4624.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 2666.6667
$ws.Range("J8").Value = 2666.6667
$ws.Range("L8").Value = 2666.6667
$ws.Range("N8").Value = -2946.6667
$ws.Range("H55").Value = 16666.666
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 16666.666
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 16666.666
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -17296.666
$ws.Range("H58").Value = 6395.7607
$ws.Range("J58").Value = 21614.334
$ws.Range("L58").Value = 21614.334
$ws.Range("N58").Value = -22020.334
$ws.Range("H99").Value = 35997.668
$ws.Range("I99").Value = 4980
$ws.Range("J99").Value = 51506.5
$ws.Range("K99").Value = 4980
$ws.Range("L99").Value = 51506.5
$ws.Range("M99").Value = -3482
$ws.Range("N99").Value = -54502.5
$ws.Range("H126").Value = 35997.668
$ws.Range("I126").Value = 4980
$ws.Range("J126").Value = 51506.5
$ws.Range("K126").Value = 14940
$ws.Range("L126").Value = 154519.5
$ws.Range("M126").Value = -12470
$ws.Range("N126").Value = -159459.5
$ws.Range("H136").Value = 6395.7607
$ws.Range("J136").Value = 21614.334
$ws.Range("L136").Value = 64843.00199999999
$ws.Range("N136").Value = -69943.00199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 229.96
$ws.Range("I26").Value = 69.15385000000001
$ws.Range("J26").Value = 404.16666
$ws.Range("K26").Value = 207.46155
$ws.Range("L26").Value = 1212.49998
$ws.Range("M26").Value = 80.53844999999998
$ws.Range("N26").Value = -1788.49998
$ws.Range("H34").Value = 2738.7778
$ws.Range("J34").Value = 3062.375
$ws.Range("L34").Value = 9187.125
$ws.Range("N34").Value = -9355.125
$ws.Range("H113").Value = 601.84
$ws.Range("I113").Value = 592
$ws.Range("J113").Value = 607.375
$ws.Range("K113").Value = 1776
$ws.Range("L113").Value = 1822.125
$ws.Range("M113").Value = 394
$ws.Range("N113").Value = -6162.125
$ws.Range("H131").Value = 777.87
$ws.Range("J131").Value = 844.8605
$ws.Range("L131").Value = 2534.5815
$ws.Range("N131").Value = -12614.5815

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 15000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 15000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 15000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -16372
$ws.Range("H65").Value = 15000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 15000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 45000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -51864
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 2446.7917
$ws.Range("I132").Value = 1760.7333
$ws.Range("K132").Value = 5282.199900000001
$ws.Range("M132").Value = -2752.199900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 6968.5
$ws.Range("J41").Value = 6968.5
$ws.Range("L41").Value = 6968.5
$ws.Range("N41").Value = -7844.5
$ws.Range("H46").Value = 2710.9167
$ws.Range("I46").Value = 510
$ws.Range("J46").Value = 3444.5557
$ws.Range("K46").Value = 510
$ws.Range("L46").Value = 3444.5557
$ws.Range("M46").Value = -322
$ws.Range("N46").Value = -3820.5557
$ws.Range("H60").Value = 25000
$ws.Range("J60").Value = 25000
$ws.Range("L60").Value = 25000
$ws.Range("N60").Value = -26018
$ws.Range("H74").Value = 15874.625
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 16713.857
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 16713.857
$ws.Range("M74").Value = -9002
$ws.Range("N74").Value = -18709.857
$ws.Range("H77").Value = 15874.625
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 16713.857
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 50141.571
$ws.Range("M77").Value = -25008
$ws.Range("N77").Value = -60125.571
$ws.Range("H132").Value = 3601.125
$ws.Range("I132").Value = 3521.2
$ws.Range("K132").Value = 10563.6
$ws.Range("M132").Value = -8033.599999999999
$ws.Range("H136").Value = 1600.2
$ws.Range("I136").Value = 1508.6666
$ws.Range("J136").Value = 1737.5
$ws.Range("K136").Value = 4525.9998
$ws.Range("L136").Value = 5212.5
$ws.Range("M136").Value = -1975.9998
$ws.Range("N136").Value = -10312.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3499.2273
$ws.Range("I132").Value = 3486.625
$ws.Range("J132").Value = 3532.8333
$ws.Range("K132").Value = 10459.875
$ws.Range("L132").Value = 10598.4999
$ws.Range("M132").Value = -7929.875
